$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition): bump 想去人数 (F column) counts ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 7506
$wsExpo.Range("F4").Value = 3562
$wsExpo.Range("F10").Value = 116
$wsExpo.Range("F12").Value = 521
$wsExpo.Range("F14").Value = 164
$wsExpo.Range("F17").Value = 357
$wsExpo.Range("F18").Value = 4229
$wsExpo.Range("F19").Value = 4229
$wsExpo.Range("F22").Value = 1036
$wsExpo.Range("F24").Value = 1919
$wsExpo.Range("F27").Value = 77
$wsExpo.Range("F28").Value = 3097
$wsExpo.Range("F29").Value = 2354
$wsExpo.Range("F34").Value = 131
$wsExpo.Range("F38").Value = 4455
$wsExpo.Range("F39").Value = 519
$wsExpo.Range("F43").Value = 851
$wsExpo.Range("F46").Value = 1680
$wsExpo.Range("F47").Value = 268
$wsExpo.Range("F50").Value = 735

# --- Sheet "全部类型" (All types): bump 想去人数 (F column) counts ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7506
$wsAll.Range("F5").Value = 3562
$wsAll.Range("F10").Value = 116
$wsAll.Range("F13").Value = 521
$wsAll.Range("F15").Value = 164
$wsAll.Range("F17").Value = 357
$wsAll.Range("F18").Value = 4229
$wsAll.Range("F19").Value = 4229
$wsAll.Range("F20").Value = 30
$wsAll.Range("F24").Value = 1036
$wsAll.Range("F26").Value = 1919
$wsAll.Range("F29").Value = 3097
$wsAll.Range("F30").Value = 2354
$wsAll.Range("F35").Value = 131
$wsAll.Range("F39").Value = 4455
$wsAll.Range("F41").Value = 519
$wsAll.Range("F45").Value = 852
$wsAll.Range("F47").Value = 1680
$wsAll.Range("F48").Value = 268
$wsAll.Range("F50").Value = 735

# --- Sheet "演出" (Performance): insert 3 new events, shifting rows, plus two F-count bumps ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("A6").Value = 5
$wsShow.Range("B6").Value = "2024-06-09"
$wsShow.Range("C6").Value = "北京·【超值5折】治愈系限定“菊次郎的夏天”·久石让&宫崎骏 主题音乐会"
$wsShow.Range("D6").Value = "朝阳北路常营陆港城20号院1号楼 常营·爱乐汇艺术空间(长楹天街店)"
$wsShow.Range("E6").Value = "2024.06.09 15:30-06.22 17:00"
$wsShow.Range("F6").Value = 0
$wsShow.Range("G6").Value = 70
$wsShow.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=86937"
$wsShow.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202406/gpDYd7CO1717563469685.jpeg"
$wsShow.Range("A7").Value = 6
$wsShow.Range("B7").Value = "2024-06-09"
$wsShow.Range("C7").Value = "北京·集结 - 超级世代！ACGN 音乐节 一周年特别庆典"
$wsShow.Range("D7").Value = "日坛北路17号日坛公园北门对面 METAL BOX"
$wsShow.Range("E7").Value = "2024.06.09 14:00-06.09 21:00"
$wsShow.Range("F7").Value = 67
$wsShow.Range("G7").Value = 80
$wsShow.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=85135"
$wsShow.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202405/ms2GHvFg1715828016187.jpeg"
$wsShow.Range("A8").Value = 7
$wsShow.Range("B8").Value = "2024-06-21"
$wsShow.Range("C8").Value = "北京·奇迹の闪耀 「UP!」巡回动漫演唱会"
$wsShow.Range("D8").Value = "亮马桥路40号(近好运街) 北京世纪剧院"
$wsShow.Range("E8").Value = "2024.06.21 19:30-06.21 21:30"
$wsShow.Range("F8").Value = 63
$wsShow.Range("G8").Value = 72
$wsShow.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=83486"
$wsShow.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202403/XKt2DiVQ1711619698950.jpeg"
$wsShow.Range("A9").Value = 8
$wsShow.Range("B9").Value = "2024-06-21"
$wsShow.Range("C9").Value = "北京·奥斯卡·罗曼耶卓（O叔）钢琴独奏音乐会"
$wsShow.Range("D9").Value = "北新华街1号 北京音乐厅"
$wsShow.Range("E9").Value = "2024.06.21 19:30-06.21 21:30"
$wsShow.Range("F9").Value = 106
$wsShow.Range("G9").Value = 480
$wsShow.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84201"
$wsShow.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202404/fAvUihAL1712887177724.jpeg"
$wsShow.Range("A10").Value = 9
$wsShow.Range("B10").Value = "2024-06-23"
$wsShow.Range("C10").Value = "北京·2024 JO☆STARS Ft. 长谷川大祐巡回演唱会"
$wsShow.Range("D10").Value = "奥园西路1号院5号楼1层2-104 福浪Live House"
$wsShow.Range("E10").Value = "2024.06.23 18:30-06.23 20:30"
$wsShow.Range("F10").Value = 30
$wsShow.Range("G10").Value = 480
$wsShow.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=85798"
$wsShow.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202405/AT8KoL6T1715915676890.png"
$wsShow.Range("A11").Value = 10
$wsShow.Range("B11").Value = "2024-06-23"
$wsShow.Range("C11").Value = "北京·仲夏绮罗日 Anisong Live Party"
$wsShow.Range("D11").Value = "大江胡同121号2幢负1层 北京门空间 TheDoorLiveHouse"
$wsShow.Range("E11").Value = "2024.06.23 13:00-06.23 16:30"
$wsShow.Range("F11").Value = 47
$wsShow.Range("G11").Value = 78
$wsShow.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=85364"
$wsShow.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202405/PFM2Be6V1715240437688.jpeg"
$wsShow.Range("A12").Value = 11
$wsShow.Range("B12").Value = "2024-06-28"
$wsShow.Range("C12").Value = "北京·“梁祝”传世经典中外小提琴名曲音乐会"
$wsShow.Range("D12").Value = "复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$wsShow.Range("E12").Value = "2024.06.28 19:30-06.28 21:00"
$wsShow.Range("F12").Value = 0
$wsShow.Range("G12").Value = 100
$wsShow.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=86906"
$wsShow.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202405/yWh8Ye2U1716537097541.png"
$wsShow.Range("A13").Value = 12
$wsShow.Range("B13").Value = "2024-06-28"
$wsShow.Range("C13").Value = "北京·《国风大赏》大型国潮音乐会×郑州歌舞剧院《唐宫夜宴》"
$wsShow.Range("D13").Value = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$wsShow.Range("E13").Value = "2024.06.28 19:30-06.28 21:00"
$wsShow.Range("F13").Value = 58
$wsShow.Range("G13").Value = 162
$wsShow.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=82587"
$wsShow.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202403/VZcJ2SJ51709882503997.jpeg"
$wsShow.Range("A14").Value = 13
$wsShow.Range("B14").Value = "2024-07-21"
$wsShow.Range("C14").Value = "北京·世界名团首席系列—— 布达佩斯节日管弦乐团弦乐四重奏音乐会"
$wsShow.Range("D14").Value = "复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$wsShow.Range("E14").Value = "2024.07.21 19:30-07.21 21:00"
$wsShow.Range("F14").Value = 1
$wsShow.Range("G14").Value = 196
$wsShow.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86891"
$wsShow.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202405/wnG2Jyvg1717049167800.png"
$wsShow.Range("A15").Value = 14
$wsShow.Range("B15").Value = "2024-07-22"
$wsShow.Range("C15").Value = "北京·石川绫子小提琴动漫音乐会"
$wsShow.Range("D15").Value = "中关村南大街33号中国国家图书馆内 国图艺术中心"
$wsShow.Range("E15").Value = "2024.07.22 19:30-07.22 21:00"
$wsShow.Range("F15").Value = 109
$wsShow.Range("G15").Value = 180
$wsShow.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=83973"
$wsShow.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202404/HhY3CS7t1712652128640.jpeg"
$wsShow.Range("A16").Value = 15
$wsShow.Range("B16").Value = "2024-07-23"
$wsShow.Range("C16").Value = "北京·巴西浪漫风情——手风琴大满贯音乐家道格拉斯·博尔萨蒂专场音乐会"
$wsShow.Range("D16").Value = "复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$wsShow.Range("E16").Value = "2024.07.23 19:30-07.23 21:00"
$wsShow.Range("F16").Value = 0
$wsShow.Range("G16").Value = 140
$wsShow.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=86922"
$wsShow.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202405/i14RABlz1716527544509.jpeg"
$wsShow.Range("A17").Value = 16
$wsShow.Range("B17").Value = "2024-07-27"
$wsShow.Range("C17").Value = "北京·“童年时光机”——《哆啦A梦》、《灌篮高手》、《狮子王》致敬童年经典动漫交响音乐会"
$wsShow.Range("D17").Value = "中关村南大街33号国家图书馆北门 国图艺术中心音乐厅"
$wsShow.Range("E17").Value = "2024.07.27 19:30-07.27 21:00"
$wsShow.Range("F17").Value = 3
$wsShow.Range("G17").Value = "不可售"
$wsShow.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=85671"
$wsShow.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202405/KV93ax2g1715669330587.jpeg"
$wsShow.Range("A18").Value = 17
$wsShow.Range("B18").Value = "2024-07-27"
$wsShow.Range("C18").Value = "北京·缤纷国图2024暑期儿童演出季 幽默钢琴莫扎特——古典音乐启蒙钢琴名曲趣味视听音乐会"
$wsShow.Range("D18").Value = "中关村南大街33号中国国家图书馆内 国图艺术中心"
$wsShow.Range("E18").Value = "2024.07.27 10:30-07.27 12:00"
$wsShow.Range("F18").Value = 1
$wsShow.Range("G18").Value = 98
$wsShow.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=86343"
$wsShow.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202405/GjHzdWRc1716782684506.jpeg"
$wsShow.Range("A19").Value = 18
$wsShow.Range("B19").Value = "2024-08-09"
$wsShow.Range("C19").Value = "北京·井草圣二 2024《夏日独白》指弹吉他音乐会"
$wsShow.Range("D19").Value = "西坝河南里2号香河园地区文化中心 多维剧场"
$wsShow.Range("E19").Value = "2024.08.09 20:00-08.09 21:30"
$wsShow.Range("F19").Value = 0
$wsShow.Range("G19").Value = 260
$wsShow.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=86938"
$wsShow.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202406/9pCUm5Pf1717642925271.jpeg"
$wsShow.Range("A20").Value = 19
$wsShow.Range("B20").Value = "2024-08-09"
$wsShow.Range("C20").Value = "北京·燃爆DNA——日本动漫原声金曲超燃演唱会"
$wsShow.Range("D20").Value = "复兴门内大街49号 民族宫大剧院"
$wsShow.Range("E20").Value = "2024.08.09 19:30-08.09 21:30"
$wsShow.Range("F20").Value = 29
$wsShow.Range("G20").Value = 180
$wsShow.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=85334"
$wsShow.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202405/WpZshtXD1715052832157.jpeg"
$wsShow.Range("A21").Value = 20
$wsShow.Range("B21").Value = "2024-08-09"
$wsShow.Range("C21").Value = "北京·阿根廷《Las Hermanas Caronni 卡洛妮姐妹二重奏》"
$wsShow.Range("D21").Value = "复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$wsShow.Range("E21").Value = "2024.08.09 19:30-08.09 21:00"
$wsShow.Range("F21").Value = 0
$wsShow.Range("G21").Value = 90
$wsShow.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=86936"
$wsShow.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202405/gdY2LOTq1716809634575.jpeg"
$wsShow.Range("A22").Value = 21
$wsShow.Range("B22").Value = "2024-08-11"
$wsShow.Range("C22").Value = "北京·Marcin Patrzalek 2024 《原声之龙》指弹吉他音乐会"
$wsShow.Range("D22").Value = "西坝河南里2号香河园地区文化中心 多维剧场"
$wsShow.Range("E22").Value = "2024.08.11 20:00-08.11 21:30"
$wsShow.Range("F22").Value = 620
$wsShow.Range("G22").Value = "不可售"
$wsShow.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=86309"
$wsShow.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202405/MEqm9GHU1716777275477.jpeg"
$wsShow.Range("A23").Value = 22
$wsShow.Range("B23").Value = "2024-08-24"
$wsShow.Range("C23").Value = "北京·最后的莫西干人——亚历桑德罗&丛林回响乐队印第安音乐品鉴会"
$wsShow.Range("D23").Value = "亮马桥路40号(近好运街) 北京世纪剧院"
$wsShow.Range("E23").Value = "2024.08.24 19:30-08.24 21:00"
$wsShow.Range("F23").Value = 3
$wsShow.Range("G23").Value = 144
$wsShow.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=86217"
$wsShow.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202405/BDyblKrJ1716427731729.jpeg"
$wsShow.Range("A24").Value = 23
$wsShow.Range("B24").Value = "2024-10-10"
$wsShow.Range("C24").Value = "北京·黑白键上的音乐地图——孩子们的钢琴协奏曲之夜"
$wsShow.Range("D24").Value = "北新华街1号 北京音乐厅"
$wsShow.Range("E24").Value = "2024.10.10 19:30-10.10 21:00"
$wsShow.Range("F24").Value = 0
$wsShow.Range("G24").Value = 144
$wsShow.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=86881"
$wsShow.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202406/K3oihoH91717474488019.jpeg"
